# DossierPlanification.docx edits
# "Fin seance mercredi 23 octobre / Modif app pas fini"

$d = $word.ActiveDocument

function Wrap-Xml($innerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-RangeWithXml($rng, $innerXml) {
    $rng.InsertXML((Wrap-Xml $innerXml))
}

function Get-CellEndPos($table, $rowIdx, $colIdx) {
    # Cell.Range.End is unreliable for multi-paragraph cells in this
    # engine (it reports the end of the first paragraph only), so derive
    # the true end boundary from the start of the following cell instead.
    $row = $table.Rows.Item($rowIdx)
    if ($colIdx -lt $row.Cells.Count) {
        return $row.Cells.Item($colIdx + 1).Range.Start - 1
    }
    if ($rowIdx -lt $table.Rows.Count) {
        $nextRow = $table.Rows.Item($rowIdx + 1)
        return $nextRow.Cells.Item(1).Range.Start - 1
    }
    return $table.Range.End - 1
}

function Replace-CellWithXml($table, $rowIdx, $colIdx, $innerXml) {
    $cell = $table.Rows.Item($rowIdx).Cells.Item($colIdx)
    $start = $cell.Range.Start
    $end = Get-CellEndPos $table $rowIdx $colIdx
    $target = $d.Range($start, $end)
    Replace-RangeWithXml $target $innerXml
}

function Replace-TextInRangeWithXml($rng, $searchText, $innerXml) {
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $searchText"
        return $false
    }
    Replace-RangeWithXml $rng $innerXml
    return $true
}

$t1 = $d.Tables.Item(1)

# --- Hunk A: Row3/C2 "Dossier Organisationnel & GANTT" -> "Rédiger d" + "ossier Organisationnel & GANTT" ---
Replace-CellWithXml $t1 3 2 '<w:body><w:p><w:r><w:t>Rédiger d</w:t></w:r><w:r><w:t>ossier Organisationnel &amp; GANTT</w:t></w:r></w:p></w:body>'
Write-Output "done A"

# --- Hunk B: Row4 "Dossier de Fabrication" ---
# C2 (Taches): "Dossier de Fabrication" -> "Rédiger d" + "ossier de Fabrication"
Replace-CellWithXml $t1 4 2 '<w:body><w:p><w:r><w:t>Rédiger d</w:t></w:r><w:r><w:t>ossier de Fabrication</w:t></w:r></w:p></w:body>'
# C3 (Livrables): "[à venir…]" -> "Dossier de Fabrication"
Replace-CellWithXml $t1 4 3 '<w:body><w:p><w:r><w:t>Dossier de Fabrication</w:t></w:r></w:p></w:body>'
Write-Output "done B"

# --- Hunk C: Row5 "Rapport de test" ---
# C2: "Rapport de test" -> "Rédiger r" + "apport de test"
Replace-CellWithXml $t1 5 2 '<w:body><w:p><w:r><w:t>Rédiger r</w:t></w:r><w:r><w:t>apport de test</w:t></w:r></w:p></w:body>'
# C3: "[à venir…]" -> "Rapport de test"
Replace-CellWithXml $t1 5 3 '<w:body><w:p><w:r><w:t>Rapport de test</w:t></w:r></w:p></w:body>'
Write-Output "done C"

# --- Hunk D: Row8/C2 "Élaboration code connexion MQTT" -> "Élaboration code connexion " + "téléphone" ---
Replace-CellWithXml $t1 8 2 '<w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">Élaboration code connexion </w:t></w:r><w:r><w:t>téléphone</w:t></w:r></w:p></w:body>'
Write-Output "done D"

# --- Hunk E: Row9/C2 "Élaboration code gestion Neopixel" -> single run (merge) ---
Replace-CellWithXml $t1 9 2 '<w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Élaboration code gestion Neopixel</w:t></w:r></w:p></w:body>'
Write-Output "done E"

# --- Hunk F: Row10/C2 "Mise en place broker MQTT" + "(raspi 0 ou directement sur ESP32)" ---
Replace-CellWithXml $t1 10 2 ('<w:body>' + `
  '<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Mise en place</w:t></w:r><w:r><w:t xml:space="preserve"> système de connexion à distance</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>(</w:t></w:r><w:r><w:t>MQTT, http,</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p>' + `
  '</w:body>')
Write-Output "done F"

Write-Output "ALL DONE"
